# Auto-generated edit script applying value updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 911
$ws.Cells.Item(5, 9).Value = 457.6
$ws.Cells.Item(5, 11).Value = 457.6
$ws.Cells.Item(5, 13).Value = -342.6
$ws.Cells.Item(9, 8).Value = 141.83333
$ws.Cells.Item(9, 9).Value = 75.2
$ws.Cells.Item(9, 10).Value = 475
$ws.Cells.Item(9, 11).Value = 75.2
$ws.Cells.Item(9, 12).Value = 475
$ws.Cells.Item(9, 13).Value = 93.8
$ws.Cells.Item(9, 14).Value = -813
$ws.Cells.Item(40, 8).Value = 4245.4814
$ws.Cells.Item(40, 9).Value = 4216.5386
$ws.Cells.Item(40, 11).Value = 4216.5386
$ws.Cells.Item(40, 13).Value = -4041.5386
$ws.Cells.Item(98, 8).Value = 2119.6365
$ws.Cells.Item(98, 10).Value = 2449
$ws.Cells.Item(98, 12).Value = 2449
$ws.Cells.Item(98, 14).Value = -5445
$ws.Cells.Item(112, 8).Value = 968.3182
$ws.Cells.Item(112, 10).Value = 973.9524
$ws.Cells.Item(112, 12).Value = 2921.8572
$ws.Cells.Item(112, 14).Value = -5137.8572
$ws.Cells.Item(122, 8).Value = 2119.6365
$ws.Cells.Item(122, 10).Value = 2449
$ws.Cells.Item(122, 12).Value = 7347
$ws.Cells.Item(122, 14).Value = -12247
$ws.Cells.Item(132, 8).Value = 7067.76
$ws.Cells.Item(132, 9).Value = 4455.15
$ws.Cells.Item(132, 11).Value = 13365.45
$ws.Cells.Item(132, 13).Value = -10835.45
$ws.Cells.Item(137, 8).Value = 4014.1628
$ws.Cells.Item(137, 10).Value = 1988.238
$ws.Cells.Item(137, 12).Value = 5964.714
$ws.Cells.Item(137, 14).Value = -11064.714
$ws.Cells.Item(139, 8).Value = 69949
$ws.Cells.Item(139, 10).Value = 69949
$ws.Cells.Item(139, 12).Value = 69949
$ws.Cells.Item(139, 14).Value = -80229

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4927.7334
$ws.Cells.Item(32, 9).Value = 4414.727
$ws.Cells.Item(32, 11).Value = 4414.727
$ws.Cells.Item(32, 13).Value = -4127.727
$ws.Cells.Item(41, 8).Value = 12731.625
$ws.Cells.Item(41, 10).Value = 24499.5
$ws.Cells.Item(41, 12).Value = 24499.5
$ws.Cells.Item(41, 14).Value = -25327.5
$ws.Cells.Item(45, 8).Value = 7466.9165
$ws.Cells.Item(45, 9).Value = 11640.077
$ws.Cells.Item(45, 11).Value = 11640.077
$ws.Cells.Item(45, 13).Value = -11263.077
$ws.Cells.Item(59, 8).Value = 47000
$ws.Cells.Item(59, 10).Value = 47000
$ws.Cells.Item(59, 12).Value = 47000
$ws.Cells.Item(59, 14).Value = -48608
$ws.Cells.Item(61, 8).Value = 4930.1553
$ws.Cells.Item(61, 9).Value = 5271.608
$ws.Cells.Item(61, 11).Value = 5271.608
$ws.Cells.Item(61, 13).Value = -5059.608
$ws.Cells.Item(136, 8).Value = 4930.1553
$ws.Cells.Item(136, 9).Value = 5271.608
$ws.Cells.Item(136, 11).Value = 15814.824
$ws.Cells.Item(136, 13).Value = -13264.824

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 33333972
$ws.Cells.Item(64, 9).Value = 41667216
$ws.Cells.Item(64, 10).Value = 1000
$ws.Cells.Item(64, 11).Value = 41667216
$ws.Cells.Item(64, 12).Value = 1000
$ws.Cells.Item(64, 13).Value = -41666991
$ws.Cells.Item(64, 14).Value = -1450
$ws.Cells.Item(67, 8).Value = 33333972
$ws.Cells.Item(67, 9).Value = 41667216
$ws.Cells.Item(67, 10).Value = 1000
$ws.Cells.Item(67, 11).Value = 41667216
$ws.Cells.Item(67, 12).Value = 1000
$ws.Cells.Item(67, 13).Value = -41666436
$ws.Cells.Item(67, 14).Value = -2560
$ws.Cells.Item(80, 8).Value = 637.9
$ws.Cells.Item(80, 9).Value = 430.625
$ws.Cells.Item(80, 11).Value = 430.625
$ws.Cells.Item(80, 13).Value = 567.375
$ws.Cells.Item(83, 8).Value = 637.9
$ws.Cells.Item(83, 9).Value = 430.625
$ws.Cells.Item(83, 11).Value = 2153.125
$ws.Cells.Item(83, 13).Value = 2838.875
$ws.Cells.Item(134, 8).Value = 3909.8948
$ws.Cells.Item(134, 9).Value = 4065.4062
$ws.Cells.Item(134, 10).Value = 3080.5
$ws.Cells.Item(134, 11).Value = 12196.2186
$ws.Cells.Item(134, 12).Value = 9241.5
$ws.Cells.Item(134, 13).Value = -9661.2186
$ws.Cells.Item(134, 14).Value = -14311.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1692.3334
$ws.Cells.Item(31, 10).Value = 2583.7
$ws.Cells.Item(31, 12).Value = 2583.7
$ws.Cells.Item(31, 14).Value = -3173.7
$ws.Cells.Item(34, 8).Value = 1692.3334
$ws.Cells.Item(34, 10).Value = 2583.7
$ws.Cells.Item(34, 12).Value = 2583.7
$ws.Cells.Item(34, 14).Value = -2987.7
$ws.Cells.Item(133, 8).Value = 154996
$ws.Cells.Item(133, 10).Value = 154996
$ws.Cells.Item(133, 12).Value = 154996
$ws.Cells.Item(133, 14).Value = -160056
$ws.Cells.Item(134, 8).Value = 5661.838
$ws.Cells.Item(134, 9).Value = 5929.3794
$ws.Cells.Item(134, 10).Value = 4692
$ws.Cells.Item(134, 11).Value = 17788.1382
$ws.Cells.Item(134, 12).Value = 14076
$ws.Cells.Item(134, 13).Value = -15253.1382
$ws.Cells.Item(134, 14).Value = -19146

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 237.83333
$ws.Cells.Item(33, 10).Value = 250.90909
$ws.Cells.Item(33, 12).Value = 1505.45454
$ws.Cells.Item(33, 14).Value = -2071.45454
$ws.Cells.Item(68, 8).Value = 2291
$ws.Cells.Item(68, 10).Value = 2649.7
$ws.Cells.Item(68, 12).Value = 7949.099999999999
$ws.Cells.Item(68, 14).Value = -9571.099999999999
$ws.Cells.Item(71, 8).Value = 2291
$ws.Cells.Item(71, 10).Value = 2649.7
$ws.Cells.Item(71, 12).Value = 23847.3
$ws.Cells.Item(71, 14).Value = -31959.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 14242.714
$ws.Cells.Item(55, 9).Value = 14633.333
$ws.Cells.Item(55, 10).Value = 13949.75
$ws.Cells.Item(55, 11).Value = 14633.333
$ws.Cells.Item(55, 12).Value = 13949.75
$ws.Cells.Item(55, 13).Value = -14306.333
$ws.Cells.Item(55, 14).Value = -14603.75
$ws.Cells.Item(58, 8).Value = 52977.168
$ws.Cells.Item(58, 10).Value = 52977.168
$ws.Cells.Item(58, 12).Value = 52977.168
$ws.Cells.Item(58, 14).Value = -53531.168
$ws.Cells.Item(70, 8).Value = 6247.0557
$ws.Cells.Item(70, 9).Value = 5934.3335
$ws.Cells.Item(70, 10).Value = 6559.778
$ws.Cells.Item(70, 11).Value = 5934.3335
$ws.Cells.Item(70, 12).Value = 6559.778
$ws.Cells.Item(70, 13).Value = -5664.3335
$ws.Cells.Item(70, 14).Value = -7099.778
$ws.Cells.Item(73, 8).Value = 6247.0557
$ws.Cells.Item(73, 9).Value = 5934.3335
$ws.Cells.Item(73, 10).Value = 6559.778
$ws.Cells.Item(73, 11).Value = 5934.3335
$ws.Cells.Item(73, 12).Value = 6559.778
$ws.Cells.Item(73, 13).Value = -4998.3335
$ws.Cells.Item(73, 14).Value = -8431.778
$ws.Cells.Item(95, 8).Value = 27549.857
$ws.Cells.Item(95, 10).Value = 27549.857
$ws.Cells.Item(95, 12).Value = 27549.857
$ws.Cells.Item(95, 14).Value = -33041.857
$ws.Cells.Item(99, 8).Value = 10907.077
$ws.Cells.Item(99, 9).Value = 6644.1113
$ws.Cells.Item(99, 10).Value = 20498.75
$ws.Cells.Item(99, 11).Value = 6644.1113
$ws.Cells.Item(99, 12).Value = 20498.75
$ws.Cells.Item(99, 13).Value = -4398.1113
$ws.Cells.Item(99, 14).Value = -24990.75
$ws.Cells.Item(126, 8).Value = 7367.8076
$ws.Cells.Item(126, 9).Value = 5967.778
$ws.Cells.Item(126, 10).Value = 8109
$ws.Cells.Item(126, 11).Value = 17903.334
$ws.Cells.Item(126, 12).Value = 24327
$ws.Cells.Item(126, 13).Value = -15433.334
$ws.Cells.Item(126, 14).Value = -29267
$ws.Cells.Item(132, 8).Value = 7370.5586
$ws.Cells.Item(132, 9).Value = 8485.521000000001
$ws.Cells.Item(132, 10).Value = 5039.273
$ws.Cells.Item(132, 11).Value = 25456.563
$ws.Cells.Item(132, 12).Value = 15117.819
$ws.Cells.Item(132, 13).Value = -22926.563
$ws.Cells.Item(132, 14).Value = -20177.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 12504412
$ws.Cells.Item(16, 9).Value = 15628278
$ws.Cells.Item(16, 10).Value = 8949.5
$ws.Cells.Item(16, 11).Value = 15628278
$ws.Cells.Item(16, 12).Value = 8949.5
$ws.Cells.Item(16, 13).Value = -15628108
$ws.Cells.Item(16, 14).Value = -9289.5
$ws.Cells.Item(22, 8).Value = 855.2
$ws.Cells.Item(22, 9).Value = 796.8570999999999
$ws.Cells.Item(22, 10).Value = 991.3333
$ws.Cells.Item(22, 11).Value = 796.8570999999999
$ws.Cells.Item(22, 12).Value = 991.3333
$ws.Cells.Item(22, 13).Value = -501.8570999999999
$ws.Cells.Item(22, 14).Value = -1581.3333
$ws.Cells.Item(27, 8).Value = 855.2
$ws.Cells.Item(27, 9).Value = 796.8570999999999
$ws.Cells.Item(27, 10).Value = 991.3333
$ws.Cells.Item(27, 11).Value = 796.8570999999999
$ws.Cells.Item(27, 12).Value = 991.3333
$ws.Cells.Item(27, 13).Value = -689.8570999999999
$ws.Cells.Item(27, 14).Value = -1205.3333
$ws.Cells.Item(46, 8).Value = 3667.5789
$ws.Cells.Item(46, 9).Value = 2119.4
$ws.Cells.Item(46, 11).Value = 2119.4
$ws.Cells.Item(46, 13).Value = -1931.4
$ws.Cells.Item(68, 8).Value = 17547298
$ws.Cells.Item(68, 9).Value = 18519370
$ws.Cells.Item(68, 11).Value = 18519370
$ws.Cells.Item(68, 13).Value = -18518621
$ws.Cells.Item(71, 8).Value = 17547298
$ws.Cells.Item(71, 9).Value = 18519370
$ws.Cells.Item(71, 11).Value = 92596850
$ws.Cells.Item(71, 13).Value = -92593106
$ws.Cells.Item(132, 8).Value = 23883.824
$ws.Cells.Item(132, 9).Value = 31274.135
$ws.Cells.Item(132, 10).Value = 4352.2856
$ws.Cells.Item(132, 11).Value = 93822.405
$ws.Cells.Item(132, 12).Value = 13056.8568
$ws.Cells.Item(132, 13).Value = -91292.405
$ws.Cells.Item(132, 14).Value = -18116.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 6998.5
$ws.Cells.Item(2, 9).Value = 6998.5
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 6998.5
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = ""
$ws.Cells.Item(2, 14).Value = -6886.5
$ws.Cells.Item(33, 8).Value = 8500
$ws.Cells.Item(33, 9).Value = 7000
$ws.Cells.Item(33, 10).Value = 10000
$ws.Cells.Item(33, 11).Value = 7000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = -6750
$ws.Cells.Item(33, 14).Value = -10500
$ws.Cells.Item(36, 8).Value = 8500
$ws.Cells.Item(36, 9).Value = 7000
$ws.Cells.Item(36, 10).Value = 10000
$ws.Cells.Item(36, 11).Value = 7000
$ws.Cells.Item(36, 12).Value = 10000
$ws.Cells.Item(36, 13).Value = -6750
$ws.Cells.Item(36, 14).Value = -10500
$ws.Cells.Item(113, 8).Value = 467.1579
$ws.Cells.Item(113, 9).Value = 430.57144
$ws.Cells.Item(113, 10).Value = 569.6
$ws.Cells.Item(113, 11).Value = 1291.71432
$ws.Cells.Item(113, 12).Value = 1708.8
$ws.Cells.Item(113, 13).Value = 878.28568
$ws.Cells.Item(113, 14).Value = -6048.8
$ws.Cells.Item(132, 8).Value = 3446.75
$ws.Cells.Item(132, 9).Value = 2558.238
$ws.Cells.Item(132, 10).Value = 9666.333000000001
$ws.Cells.Item(132, 11).Value = 7674.714
$ws.Cells.Item(132, 12).Value = 28998.999
$ws.Cells.Item(132, 13).Value = -5144.714
$ws.Cells.Item(132, 14).Value = -34058.999
$ws.Cells.Item(136, 8).Value = 2257.9443
$ws.Cells.Item(136, 9).Value = 2208.0645
$ws.Cells.Item(136, 10).Value = 2567.2
$ws.Cells.Item(136, 11).Value = 6624.193499999999
$ws.Cells.Item(136, 12).Value = 7701.599999999999
$ws.Cells.Item(136, 13).Value = -4074.193499999999
$ws.Cells.Item(136, 14).Value = -12801.6
